$wb = $excel.ActiveWorkbook

$wsCap = $wb.Worksheets.Item("Capacity")
$wsRem = $wb.Worksheets.Item("Remove_units_by_node")

# ---------------------------------------------------------------------------
# Capacity sheet: add the "Distributed Energy" / 2030 rows, mirroring the
# existing "National Trends"/2025 (rows 2-3,6) and "Distributed Energy"/2040
# (rows 4-5,7) blocks for row 8 (Lower reservoir turbine), row 9 (Upper
# reservoir turbine) and row 10 (Lower reservoir turbine peaker).
# ---------------------------------------------------------------------------

# Row 8: copy formatting from row 4 (same pattern, just different year/scenario)
$wsCap.Range("A4:J4").Copy()
$wsCap.Range("A8:J8").PasteSpecial(-4122)
$wsCap.Range("I8:J8").Clear()
$wsCap.Range("A8").Value = "FI00"
$wsCap.Range("C8").Value = "Lower reservoir turbine"
$wsCap.Range("D8").Value = "Distributed Energy"
$wsCap.Range("E8").Value = 2030
$wsCap.Range("F8").Value = 2400

# Row 9: copy formatting from row 5
$wsCap.Range("A5:J5").Copy()
$wsCap.Range("A9:J9").PasteSpecial(-4122)
$wsCap.Range("I9:J9").Clear()
$wsCap.Range("A9").Value = "FI00"
$wsCap.Range("C9").Value = "Upper reservoir turbine"
$wsCap.Range("D9").Value = "Distributed Energy"
$wsCap.Range("E9").Value = 2030
$wsCap.Range("F9").ClearContents()
$wsCap.Range("H9").Value = 1800

# Row 10: copy formatting from row 7 (non-contiguous cells: A, C, D, E, F only)
$wsCap.Range("A10:J10").Clear()
$wsCap.Range("A7").Copy()
$wsCap.Range("A10").PasteSpecial(-4122)
$wsCap.Range("C7").Copy()
$wsCap.Range("C10").PasteSpecial(-4122)
$wsCap.Range("D7").Copy()
$wsCap.Range("D10").PasteSpecial(-4122)
$wsCap.Range("E7").Copy()
$wsCap.Range("E10").PasteSpecial(-4122)
$wsCap.Range("F7").Copy()
$wsCap.Range("F10").PasteSpecial(-4122)
$wsCap.Range("A10").Value = "FI00"
$wsCap.Range("C10").Value = "Lower reservoir turbine peaker"
$wsCap.Range("D10").Value = "Distributed Energy"
$wsCap.Range("E10").Value = 2030
$wsCap.Range("F10").Value = 200

# Rows 2, 3 and 6 become hidden (outline / manual row hide)
$wsCap.Rows.Item(2).Hidden = $true
$wsCap.Rows.Item(3).Hidden = $true
$wsCap.Rows.Item(6).Hidden = $true

# Re-apply the AutoFilter, now filtering column D (Scenario) on "Distributed Energy"
$wsCap.Range("A1:J65").AutoFilter(4, @("Distributed Energy", ""), 7)

# ---------------------------------------------------------------------------
# Remove_units_by_node sheet: add the FI00 / Reservoir / Distributed Energy /
# 2030 row, mirroring rows 2-3.
# ---------------------------------------------------------------------------
$wsRem.Range("A2:D2").Copy()
$wsRem.Range("A4:D4").PasteSpecial(-4122)
$wsRem.Range("A4").Value = "FI00"
$wsRem.Range("B4").Value = "Reservoir"
$wsRem.Range("C4").Value = "Distributed Energy"
$wsRem.Range("D4").Value = 2030

# ---------------------------------------------------------------------------
# Selections / active sheet: Capacity becomes the active tab with E15
# selected; Remove_units_by_node keeps C12 selected but is no longer active.
# ---------------------------------------------------------------------------
$wsRem.Range("C12").Select()
$wsCap.Activate()
$wsCap.Range("E15").Select()

$wb.Save()
